$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Headers: wdHeaderFooterPrimary=1 (default header -> header2.xml, logo "image1.jpg"),
#          wdHeaderFooterFirstPage=2 (first-page header -> header1.xml, logo "image1.jpg")
for ($i = 1; $i -le 2; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        for ($j = 1; $j -le $hdr.Range.InlineShapes.Count; $j++) {
            $ish = $hdr.Range.InlineShapes.Item($j)
            if ($ish.AlternativeText -eq "BTec_Logo-Orange") {
                $ish.Name = "image2.jpg"
            }
        }
    }
}

# Footers: wdHeaderFooterPrimary=1 (default footer -> footer2.xml, logo "image2.png"),
#          wdHeaderFooterFirstPage=2 (first-page footer -> footer1.xml, logo "image2.png")
for ($i = 1; $i -le 2; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        for ($j = 1; $j -le $ftr.Range.InlineShapes.Count; $j++) {
            $ish = $ftr.Range.InlineShapes.Item($j)
            if ($ish.AlternativeText -like "*PearsonLogo*") {
                $ish.Name = "image1.png"
            }
        }
    }
}
